# add unjoin_room & fix join_room
# ------------------------------------------------------------------
# This script reproduces, via Excel COM automation, the hand-edit that:
#   1) removes the blank separator row that used to sit above the
#      "room_join" (방가입) block, and appends a new FAIL row
#      (NOT_FOUND_ROOM / 잘못된 방 입력) to the bottom of that block,
#   2) removes the (now redundant) blank separator rows that used to
#      sit above "유저 팔로우", "댓글 가져오기", "팔로우 중인지 확인"
#      and "유저 언 팔로우",
#   3) appends a brand new "방탈퇴" (room_unjoin_check) block at the
#      bottom of the sheet, complete with merged/styled header cell
#      and a hyperlink pointing at the new endpoint, and
#   4) leaves the selection on I50 (matching the author's last
#      recorded cursor position).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("php")

# --- 1) room_join (방가입) block -----------------------------------
# Before: row38 = blank separator, row39-42 = 방가입 block (header +
#         3 rows), row43 = next block header ("게시글가져오기").
# Delete the blank separator above the block so the header moves up
# to row 38.
$ws.Rows.Item(38).Delete()

# The block now occupies rows 38-41 (row41 = MEM_LIMIT / 인원수 초과).
# Insert a blank row right after it (pushing "게시글가져오기" back down
# to row43 where it already was) and fill it with the new FAIL case.
$ws.Rows.Item(42).Insert()
$ws.Range("J42").Value = "NOT_FOUND_ROOM"
$ws.Range("K42").Value = "잘못된 방 입력"

# --- 2) drop the blank separator rows that are no longer present ---
# After step 1 these separators live at rows 47, 52, 56 and 61 (they
# shift up by one apiece as each prior one is removed).
$ws.Rows.Item(47).Delete()
$ws.Rows.Item(51).Delete()
$ws.Rows.Item(54).Delete()
$ws.Rows.Item(58).Delete()

# --- 3) new "방탈퇴" (room_unjoin_check) block at rows 61-64 --------
$ws.Range("A61").Value = "방탈퇴"
$ws.Range("J61").Value = "SUCCESS"

$ws.Range("A62:E62").Merge()
$ws.Range("A62").Style = $ws.Range("A59").Style
$ws.Range("A62").Value = "http://lemontree.dothome.co.kr/pinbox/room/room_unjoin_check"
$ws.Range("F62").Value = "user_token"
$ws.Range("G62").Value = "유저 토큰"
$ws.Range("J62").Value = "FAIL"
$ws.Range("K62").Value = "성공적으로 탈퇴됨"

$ws.Range("F63").Value = "room_token"
$ws.Range("G63").Value = "방 토큰"
$ws.Range("J63").Value = "NOT_FOUND_ROOM"
$ws.Range("K63").Value = "실패"

$ws.Range("K64").Value = "잘못된 방 입력"

$ws.Hyperlinks.Add($ws.Range("A62"), "http://lemontree.dothome.co.kr/pinbox/room/room_unjoin_check")

# --- 4) restore cursor/selection ------------------------------------
$ws.Activate()
$ws.Range("I50").Select()

Write-Output "edit complete"
